$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "Output"
$ws.Range("D16").Value = "HO_VectorBox/modes"
$ws.Range("C17").Value = "Output"
$ws.Range("D17").Value = "HO_VectorBox/motorControlSlow"

$ws.Range("D7").Value = "HO_VectorBox/motorControlSlow"
$ws.Range("D10").Value = "HO_VectorBox/motorControlFast"
$ws.Range("D13").Value = "HO_VectorBox/motorControlFast"
$ws.Range("D15").Value = "HO_VectorBox/motorControlSlow"

$ws.Range("G7").Select()
